$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E4").Value = "2016-03-25 00:57:56"
$wsZh.Range("H4").Value = "2016-03-25 00:58:26"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E4").Value = "2016-03-25 00:58:01"
$wsDe.Range("H4").Value = "2016-03-25 00:58:33"
